# Rename three slide layouts and strip their now-unused custom
# placeholder shapes, reverting them back to the stock Office layout
# names/contents.

$p = $ppt.ActivePresentation
$layouts = $p.SlideMaster.CustomLayouts

# --- slideLayout5.xml: "Strategy" -> "Comparison" -------------------
$layout5 = $layouts.Item(5)
$layout5.Name = "Comparison"
$layout5.Shapes.Item("Strat Title").Delete()
$layout5.Shapes.Item("Management").Delete()
$layout5.Shapes.Item("Segments Pie").Delete()
$layout5.Shapes.Item("News").Delete()

# --- slideLayout6.xml: "Financials" -> "Title Only" ------------------
$layout6 = $layouts.Item(6)
$layout6.Name = "Title Only"
$layout6.Shapes.Item("Fin Title").Delete()
$layout6.Shapes.Item("Fin Table").Delete()
$layout6.Shapes.Item("EBITDA Chart").Delete()
$layout6.Shapes.Item("Deal Score").Delete()

# --- slideLayout7.xml: "Executive Summary" -> "Blank" -----------------
$layout7 = $layouts.Item(7)
$layout7.Name = "Blank"
$layout7.Shapes.Item("Company Name").Delete()
$layout7.Shapes.Item("Ticker Price").Delete()
$layout7.Shapes.Item("Summary Bullets").Delete()
$layout7.Shapes.Item("Price Chart").Delete()
$layout7.Shapes.Item("Key Metrics").Delete()
